$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.333087921142578
$ws.Range("B1").Value = 2.470479249954224
$ws.Range("C1").Value = 6.017116069793701
$ws.Range("D1").Value = 1.880101084709167
$ws.Range("E1").Value = 1.269769549369812
